# Criação da tabela "Tabela de Dia da Semana" (Plan2), espelhando a
# tabela existente "Tabela de Tipos de Pessoas".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cabeçalho mesclado A5:B5, com a mesma formatação (negrito + centralizado)
# do cabeçalho A1:B1, reaproveitando o estilo existente via copiar/colar
# formatos (evita criar estilos novos na planilha).
$ws.Range("A1:B1").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5:B5").Merge()
$ws.Range("A5").Value = "Tabela de Dia da Semana"
$excel.CutCopyMode = $false

# Linhas de dados: número do dia (1-7) e nome do dia da semana.
$days = @("Segunda-feira", "Terça-feira", "Quarta-feira", "Quinta-feira", "Sexta-feira", "Sábado", "Domingo")
for ($i = 0; $i -lt $days.Length; $i++) {
    $row = 6 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $days[$i]
}

# Seleciona a tabela recém-criada (replica o estado final do autor).
$ws.Range("A6:B12").Select()
